$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new log row (row 3) under the existing header/data rows.
$ws.Range("A3").Value = "edit1"
$ws.Range("B3").Value = "riya-morankar"
$ws.Range("C3").Value = "Merged"
$ws.Range("D3").Value = "desc"

# Force the date column to stay literal text ("2025-06-17") instead of
# being auto-converted into a date serial number.
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2025-06-17"

$ws.Range("F3").Value = "N/A"
